$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 111, shifting existing rows 111:120 down to 112:121
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new record
$ws.Range("A111").Value = 11
$ws.Range("B111").Value = "Vega Monumental Concepción"
$ws.Range("C111").Value = "Bíobío"
$ws.Range("D111").Value = 44461
$ws.Range("E111").Value = 8
$ws.Range("F111").Value = 100112045
$ws.Range("G111").Value = "Zapallo"
$ws.Range("H111").Value = "Paine"
$ws.Range("I111").Value = "1a (guarda)"
$ws.Range("J111").Value = 300
$ws.Range("K111").Value = 200
$ws.Range("L111").Value = 220
$ws.Range("M111").Value = 210
$ws.Range("N111").Value = "$/kilo (volumen en unidades)"
$ws.Range("O111").Value = "Región de O'Higgins"
$ws.Range("P111").Value = 210
$ws.Range("Q111").Value = 1
$ws.Range("R111").Value = "Hortaliza"

# Ensure the number format of the new date cell matches the rest of column D
$ws.Range("D111").NumberFormat = $ws.Range("D112").NumberFormat
